$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 383.9375
$ws.Range("I9").Value = 364.54544
$ws.Range("J9").Value = 426.6
$ws.Range("K9").Value = 364.54544
$ws.Range("L9").Value = 426.6
$ws.Range("M9").Value = -195.54544
$ws.Range("N9").Value = -764.6
$ws.Range("H38").Value = 20218
$ws.Range("J38").Value = 55000
$ws.Range("L38").Value = 165000
$ws.Range("N38").Value = -165744
$ws.Range("H41").Value = 70
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").Value = ""
$ws.Range("H62").Value = 10229.333
$ws.Range("J62").Value = 12344.5
$ws.Range("L62").Value = 12344.5
$ws.Range("N62").Value = -13592.5
$ws.Range("H65").Value = 10229.333
$ws.Range("J65").Value = 12344.5
$ws.Range("L65").Value = 61722.5
$ws.Range("N65").Value = -67962.5
$ws.Range("H74").Value = 7422.6665
$ws.Range("I74").Value = 7561.2
$ws.Range("K74").Value = 7561.2
$ws.Range("M74").Value = -6625.2
$ws.Range("H77").Value = 7422.6665
$ws.Range("I77").Value = 7561.2
$ws.Range("K77").Value = 37806
$ws.Range("M77").Value = -33126
$ws.Range("H96").Value = 1908.8889
$ws.Range("I96").Value = 3062.4
$ws.Range("J96").Value = 467
$ws.Range("K96").Value = 9187.2
$ws.Range("L96").Value = 1401
$ws.Range("M96").Value = -7814.200000000001
$ws.Range("N96").Value = -4147
$ws.Range("H106").Value = 5559965
$ws.Range("I106").Value = 6671258.5
$ws.Range("K106").Value = 6671258.5
$ws.Range("M106").Value = -6670627.5
$ws.Range("H121").Value = 2648.1667
$ws.Range("J121").Value = 2648.1667
$ws.Range("L121").Value = 7944.500100000001
$ws.Range("N121").Value = -11438.5001
$ws.Range("H127").Value = 1216.3334
$ws.Range("J127").Value = 2749
$ws.Range("L127").Value = 8247
$ws.Range("N127").Value = -18167
$ws.Range("H132").Value = 3022.7354
$ws.Range("I132").Value = 2940.652
$ws.Range("K132").Value = 8821.956
$ws.Range("M132").Value = -6291.956
$ws.Range("H141").Value = 7197.8184
$ws.Range("I141").Value = 5935.476
$ws.Range("K141").Value = 17806.428
$ws.Range("M141").Value = -12626.428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1862
$ws.Range("I2").Value = 1610.25
$ws.Range("J2").Value = 2113.75
$ws.Range("K2").Value = 1610.25
$ws.Range("L2").Value = 2113.75
$ws.Range("M2").Value = -1497.25
$ws.Range("N2").Value = -2339.75
$ws.Range("H10").Value = 13338000
$ws.Range("I10").Value = 20005000
$ws.Range("J10").Value = 4000
$ws.Range("K10").Value = 20005000
$ws.Range("L10").Value = 4000
$ws.Range("M10").Value = -20004830
$ws.Range("N10").Value = -4340
$ws.Range("H32").Value = 17369038
$ws.Range("I32").Value = 17943956
$ws.Range("J32").Value = 11907319
$ws.Range("K32").Value = 17943956
$ws.Range("L32").Value = 11907319
$ws.Range("M32").Value = -17943669
$ws.Range("N32").Value = -11907893
$ws.Range("H74").Value = 2598.8462
$ws.Range("I74").Value = 2565.4583
$ws.Range("J74").Value = 2999.5
$ws.Range("K74").Value = 2565.4583
$ws.Range("L74").Value = 2999.5
$ws.Range("M74").Value = -1691.4583
$ws.Range("N74").Value = -4747.5
$ws.Range("H77").Value = 2598.8462
$ws.Range("I77").Value = 2565.4583
$ws.Range("J77").Value = 2999.5
$ws.Range("K77").Value = 12827.2915
$ws.Range("L77").Value = 14997.5
$ws.Range("M77").Value = -8459.2915
$ws.Range("N77").Value = -23733.5
$ws.Range("H98").Value = 100000
$ws.Range("J98").Value = 100000
$ws.Range("L98").Value = 100000
$ws.Range("N98").Value = -105990
$ws.Range("H116").Value = 1862
$ws.Range("I116").Value = 1610.25
$ws.Range("J116").Value = 2113.75
$ws.Range("K116").Value = 1610.25
$ws.Range("L116").Value = 2113.75
$ws.Range("M116").Value = 683.75
$ws.Range("N116").Value = -6701.75
$ws.Range("H122").Value = 8384.186
$ws.Range("I122").Value = 6447.3
$ws.Range("K122").Value = 19341.9
$ws.Range("M122").Value = -16891.9
$ws.Range("H132").Value = 3355.0356
$ws.Range("I132").Value = 2717.1052
$ws.Range("K132").Value = 8151.3156
$ws.Range("M132").Value = -5621.3156

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1862
$ws.Range("I3").Value = 1610.25
$ws.Range("J3").Value = 2113.75
$ws.Range("K3").Value = 1610.25
$ws.Range("L3").Value = 2113.75
$ws.Range("M3").Value = -1496.25
$ws.Range("N3").Value = -2341.75
$ws.Range("H86").Value = 994
$ws.Range("I86").Value = 994
$ws.Range("K86").Value = 994
$ws.Range("M86").Value = 129
$ws.Range("H89").Value = 994
$ws.Range("I89").Value = 994
$ws.Range("K89").Value = 4970
$ws.Range("M89").Value = 646
$ws.Range("H94").Value = 1096.0344
$ws.Range("I94").Value = 819.875
$ws.Range("J94").Value = 1435.9231
$ws.Range("K94").Value = 819.875
$ws.Range("L94").Value = 1435.9231
$ws.Range("M94").Value = -368.875
$ws.Range("N94").Value = -2337.9231
$ws.Range("H99").Value = 2402.4783
$ws.Range("I99").Value = 2014.2778
$ws.Range("J99").Value = 3800
$ws.Range("K99").Value = 2014.2778
$ws.Range("L99").Value = 3800
$ws.Range("M99").Value = -516.2778000000001
$ws.Range("N99").Value = -6796
$ws.Range("H134").Value = 3511941.8
$ws.Range("J134").Value = 6450
$ws.Range("L134").Value = 19350
$ws.Range("N134").Value = -24420

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 61897.25
$ws.Range("J28").Value = 61897.25
$ws.Range("L28").Value = 61897.25
$ws.Range("N28").Value = -62387.25
$ws.Range("H105").Value = 3739.2
$ws.Range("I105").Value = 2898.5
$ws.Range("J105").Value = 4299.6665
$ws.Range("K105").Value = 2898.5
$ws.Range("L105").Value = 4299.6665
$ws.Range("M105").Value = -1151.5
$ws.Range("N105").Value = -7793.6665
$ws.Range("H122").Value = 4412.864
$ws.Range("I122").Value = 3206.875
$ws.Range("J122").Value = 7628.8335
$ws.Range("K122").Value = 9620.625
$ws.Range("L122").Value = 22886.5005
$ws.Range("M122").Value = -7170.625
$ws.Range("N122").Value = -27786.5005
$ws.Range("H134").Value = 3304
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 3304
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 9912
$ws.Range("M134").Value = ""
$ws.Range("N134").Value = -14982

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 793.1111
$ws.Range("I60").Value = 885.5
$ws.Range("J60").Value = 719.2
$ws.Range("K60").Value = 2656.5
$ws.Range("L60").Value = 2157.6
$ws.Range("M60").Value = -2405.5
$ws.Range("N60").Value = -2659.6
$ws.Range("H113").Value = 4249.0835
$ws.Range("J113").Value = 4544.4546
$ws.Range("L113").Value = 13633.3638
$ws.Range("N113").Value = -17973.3638
$ws.Range("H129").Value = 1679.1666
$ws.Range("I129").Value = 452
$ws.Range("J129").Value = 2292.75
$ws.Range("K129").Value = 1356
$ws.Range("L129").Value = 6878.25
$ws.Range("M129").Value = 3644
$ws.Range("N129").Value = -16878.25
$ws.Range("H131").Value = 1800.5518
$ws.Range("J131").Value = 1832.92
$ws.Range("L131").Value = 5498.76
$ws.Range("N131").Value = -15578.76
$ws.Range("H132").Value = 1218.375
$ws.Range("J132").Value = 1591.8
$ws.Range("L132").Value = 14326.2
$ws.Range("N132").Value = -19386.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3113.4443
$ws.Range("J80").Value = 3666.3333
$ws.Range("L80").Value = 3666.3333
$ws.Range("N80").Value = -5662.3333
$ws.Range("H83").Value = 3113.4443
$ws.Range("J83").Value = 3666.3333
$ws.Range("L83").Value = 18331.6665
$ws.Range("N83").Value = -28315.6665
$ws.Range("H97").Value = 1013.4545
$ws.Range("I97").Value = 706.625
$ws.Range("J97").Value = 1831.6666
$ws.Range("K97").Value = 706.625
$ws.Range("L97").Value = 1831.6666
$ws.Range("M97").Value = -210.625
$ws.Range("N97").Value = -2823.6666
$ws.Range("H122").Value = 1115.6
$ws.Range("I122").Value = 1115.6
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3346.8
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -896.7999999999997
$ws.Range("N122").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1545.4762
$ws.Range("I93").Value = 1476.4546
$ws.Range("J93").Value = 1621.4
$ws.Range("K93").Value = 1476.4546
$ws.Range("L93").Value = 1621.4
$ws.Range("M93").Value = -228.4546
$ws.Range("N93").Value = -4117.4
$ws.Range("H122").Value = 13417.211
$ws.Range("I122").Value = 13437.12
$ws.Range("K122").Value = 40311.36
$ws.Range("M122").Value = -37861.36
$ws.Range("H136").Value = 9076.857
$ws.Range("I136").Value = 6147.6
$ws.Range("J136").Value = 16400
$ws.Range("K136").Value = 18442.8
$ws.Range("L136").Value = 49200
$ws.Range("M136").Value = -15892.8
$ws.Range("N136").Value = -54300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 25000
$ws.Range("J17").Value = 25000
$ws.Range("L17").Value = 25000
$ws.Range("N17").Value = -25344
$ws.Range("H81").Value = 4047.389
$ws.Range("I81").Value = 3841.25
$ws.Range("J81").Value = 4106.2856
$ws.Range("K81").Value = 7682.5
$ws.Range("L81").Value = 8212.5712
$ws.Range("M81").Value = -6621.5
$ws.Range("N81").Value = -10334.5712
$ws.Range("H84").Value = 4047.389
$ws.Range("I84").Value = 3841.25
$ws.Range("J84").Value = 4106.2856
$ws.Range("K84").Value = 38412.5
$ws.Range("L84").Value = 41062.856
$ws.Range("M84").Value = -33108.5
$ws.Range("N84").Value = -51670.856
$ws.Range("H107").Value = 486.55554
$ws.Range("J107").Value = 546.8333
$ws.Range("L107").Value = 1640.4999
$ws.Range("N107").Value = -5480.4999
$ws.Range("H122").Value = 5458.8696
$ws.Range("I122").Value = 5050.579
$ws.Range("K122").Value = 15151.737
$ws.Range("M122").Value = -12701.737
$ws.Range("H132").Value = 2389.2778
$ws.Range("I132").Value = 2233.8667
$ws.Range("K132").Value = 6701.6001
$ws.Range("M132").Value = -4171.6001
$ws.Range("H136").Value = 1303.8235
$ws.Range("I136").Value = 982.8571
$ws.Range("J136").Value = 2801.6667
$ws.Range("K136").Value = 2948.5713
$ws.Range("L136").Value = 8405.000100000001
$ws.Range("M136").Value = -398.5712999999996
$ws.Range("N136").Value = -13505.0001
